# GoodInfo_v2 - 2021.12.07
# Append the next day's row to the 未實現報酬率 (unrealized return) table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to be read as plain text (matching the existing
# rows, which store "2021-12-0x" as text rather than a date serial), then
# restore the default cell style so the new row doesn't pick up a stray
# number-format style.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2021-12-07"
$ws.Range("A5").Style = $ws.Range("A2").Style

$ws.Range("B5").Value = 3.71
